# Fix Excel test case for Open Menu issue
# - Clears the contents of row 5 (A5:B5), which previously held the
#   "locked_out_user" / "secret_sauce" test row, while keeping A5's
#   existing cell style.
# - Selects the full row 5 (A5:XFD5) with active cell A5, matching the
#   updated selection state saved with the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cell contents for A5:B5 (keeps existing formatting/style on A5)
$ws.Range("A5:B5").ClearContents()

# Update the sheet's selection to the full row A5:XFD5, active cell A5
$ws.Range("A5:XFD5").Select()
